$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Order matters: new shared-string values must be introduced in the same
# sequence they appear in the target workbook so the rebuilt sharedStrings
# table lines up. Re-used (pre-existing) values can be written at any point.

# Row 2 (Lucknow): LiveTemp + Wind refreshed
$ws.Range("B2").Value = "27.00"
$ws.Range("K2").Value = "3.29"

# Row 3 (Bhopal): Wind, Weather Condition refreshed; LiveTemp/Compare/Humidity updated
$ws.Range("K3").Value = "3.60"
$ws.Range("E3").Value = "Clouds and scattered clouds"
$ws.Range("B3").Value = "27.00"
$ws.Range("D3").Value = "Temperatues are NOT within Variance Range"
$ws.Range("H3").Value = "83"

# Row 4 (Ajmer): LiveTemp, Wind, Weather Condition refreshed; Compare/Humidity updated
$ws.Range("B4").Value = "29.82"
$ws.Range("K4").Value = "7.93"
$ws.Range("E4").Value = "Clouds and overcast clouds"
$ws.Range("D4").Value = "Temperatues are within Variance Range"
$ws.Range("H4").Value = "65"

# Row 5 (Coimbatore): LiveTemp refreshed; Weather Condition/Compare/Humidity/Wind updated
$ws.Range("B5").Value = "26.00"
$ws.Range("E5").Value = "Clouds and scattered clouds"
$ws.Range("G5").Value = "Both portals show slightly similar Weather conditions!"
$ws.Range("H5").Value = "83"
$ws.Range("K5").Value = "3.60"

# Row 6 (Mumbai): Wind refreshed
$ws.Range("K6").Value = "2.60"

# Row 7 (Kolkata): LiveTemp, Humidity, Wind refreshed
$ws.Range("B7").Value = "28.00"
$ws.Range("H7").Value = "94"
$ws.Range("K7").Value = "5.95"
